$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For Price-column values that look like genuine numbers (e.g. "238.58"),
# force the cell to Text format first so Excel keeps the exact original
# string (trailing zeros, precision) instead of silently parsing it into a
# floating point number. Non-numeric-looking strings (coin names, links,
# percentages, and "thousands-dot" prices like "30.465.96") are left alone
# since Excel already keeps those as plain text automatically.

$ws.Range('D2').Value = '30.465.96'
$ws.Range('E2').Value = '  -0.18%  '

$ws.Range('D3').Value = '1.899.81'
$ws.Range('E3').Value = '  +1.36%  '

$ws.Range('E4').Value = '  +0.07%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '238.58'
$ws.Range('E5').Value = '  +0.99%  '

$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$ws.Range('E6').Value = '  +0.09%  '

$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.4894'
$ws.Range('E7').Value = '  +0.57%  '

$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.2920'
$ws.Range('E8').Value = '  +0.95%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.06666'
$ws.Range('E9').Value = '  +0.02%  '

$ws.Range('D10').Value = '1.903.13'
$ws.Range('E10').Value = '  +1.68%  '

$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '16.96'
$ws.Range('E11').Value = '  +2.32%  '

$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '0.07321'
$ws.Range('E12').Value = '  +1.37%  '

$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '5.185'
$ws.Range('E13').Value = '  +3.72%  '

$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '88.07'
$ws.Range('E14').Value = '  -1.50%  '

$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.6664'
$ws.Range('E15').Value = '  +2.02%  '

$ws.Range('D16').Value = '30.446.13'
$ws.Range('E16').Value = '  -0.05%  '

$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '13.43'
$ws.Range('E17').Value = '  +3.26%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '0.000007856'
$ws.Range('E18').Value = '  +0.43%  '

$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$ws.Range('E19').Value = '  +0.10%  '

$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '5.441'
$ws.Range('E20').Value = '  +15.04%  '

$ws.Range('D21').Value = '2.140.00'
$ws.Range('E21').Value = '  +1.26%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$ws.Range('E22').Value = '  +0.14%  '

$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '195.27'
$ws.Range('E23').Value = '  -8.44%  '

$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '6.134'
$ws.Range('E24').Value = '  +0.13%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '9.473'
$ws.Range('E25').Value = '  +1.13%  '

$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '163.31'
$ws.Range('E26').Value = '  +4.54%  '

$ws.Range('E27').Value = '  -3.81%  '

$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '1.937'
$ws.Range('E28').Value = '  +5.99%  '

$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '1.473'
$ws.Range('E29').Value = '  +4.35%  '

$ws.Range('E30').Value = '  +1.60%  '

$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '0.09158'
$ws.Range('E31').Value = '  +1.32%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '4.128'
$ws.Range('E32').Value = '  +5.35%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '0.05168'
$ws.Range('E33').Value = '  +0.96%  '

$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '0.7393'
$ws.Range('E34').Value = '  +2.03%  '

$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '1.107'
$ws.Range('E35').Value = '  +2.85%  '

$ws.Range('E36').Value = '  +1.60%  '

$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '0.01845'
$ws.Range('E37').Value = '  +1.70%  '

$ws.Range('E38').Value = '  +0.68%  '

$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.9237'
$ws.Range('E39').Value = '  +0.50%  '

$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '2.062'
$ws.Range('E40').Value = '  +0.92%  '

$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.4404'
$ws.Range('E41').Value = '  -0.03%  '

$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '106.99'
$ws.Range('E42').Value = '  +2.40%  '

$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '5.897'
$ws.Range('E43').Value = '  +2.76%  '

$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '0.9951'
$ws.Range('E44').Value = '  +0.11%  '

$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '68.74'
$ws.Range('E45').Value = '  +20.57%  '

$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '0.1370'
$ws.Range('E46').Value = '  +3.32%  '

$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '7.558'
$ws.Range('E47').Value = '  +3.18%  '

$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '8.986'
$ws.Range('E48').Value = '  +3.89%  '

$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '34.93'
$ws.Range('E49').Value = '  +5.19%  '

$ws.Range('E50').Value = '  +0.06%  '

$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.3925'
$ws.Range('E51').Value = '  -2.11%  '
